$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.934
$ws.Range("C8").Value = -12.672
$ws.Range("B12").Value = 5.821
$ws.Range("C12").Value = -12.911
$ws.Range("C14").Value = -12.049
$ws.Range("C22").Value = -12.929
